# edit.ps1 — applies the "added a start of section 5 to the cheatsheet
# and a few more minor edits to MT materials" commit.
#
#   1. Bump every cached datetimeFigureOut field (master, all layouts,
#      notes master) from 12/6/2023 to 12/7/2023.
#   2. Slide 4: "R ignores whitespace" -> "R ignores whitespace between things"
#   3. Slide 43: grow the "R ignores whitespace, which includes <return>"
#      textbox and extend its wording the same way.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholders: 12/6/2023 -> 12/7/2023
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "12/6/2023") {
                $tr.Text = "12/7/2023"
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master (its date placeholder only accepts writes through the
# HeadersFooters facade, not through the shape's TextFrame directly)
$nmShape = $p.NotesMaster.Shapes.Item(2)
if ($nmShape.Name -like "Date Placeholder*" -and $nmShape.TextFrame.TextRange.Text -eq "12/6/2023") {
    $p.NotesMaster.HeadersFooters.DateAndTime.Text = "12/7/2023"
}

# ---------------------------------------------------------------------
# 2. Slide 4 - "R ignores whitespace" -> "R ignores whitespace between things"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $shp = $s4.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange
    $full = $tr.Text
    $needle = "R ignores whitespace"
    $idx = $full.IndexOf($needle)
    if ($idx -ge 0 -and $full.Substring($idx, $needle.Length + 1) -ne "R ignores whitespace ") {
        $sub = $tr.Characters($idx + 1, $needle.Length)
        $sub.Text = "R ignores whitespace between things"
    }
}

# ---------------------------------------------------------------------
# 3. Slide 43 - resize textbox + extend wording
# ---------------------------------------------------------------------
$s43 = $p.Slides.Item(43)
for ($i = 1; $i -le $s43.Shapes.Count; $i++) {
    $shp = $s43.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange
    $full = $tr.Text
    $needle = "R ignores whitespace, which includes"
    if ($full.IndexOf($needle) -ge 0) {
        # grow the box (only the height changes)
        $shp.Height = 103.58787401574803

        $idx = $full.IndexOf("R ignores whitespace")
        $sub = $tr.Characters($idx + 1, "R ignores whitespace".Length)
        $sub.Text = "R ignores whitespace between things"
    }
}
